# Corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 duplicated the header row (arg/code/codnbody/...) with no
# formatting; it was spurious and is removed entirely. Deleting it shifts
# the metric rows (old 3..9) up by one (new 2..8).
$ws.Rows(2).Delete()

# A1 used to hold the stray "Unnamed: 0" label (the index column header from
# the source dataframe); it is now blank. Clear the value before restyling
# so the (now empty) cell record is preserved rather than dropped.
$ws.Range("A1").Value = ""

# Row 1 no longer carries the bold/border/center-aligned header style - put
# it back to the workbook's default "Normal" style.
$ws.Range("A1:R1").Style = "Normal"

# The "code" (C) and "summary" (N) columns were recomputed with corrected
# pre/post/total fixation figures for every metric row.
$ws.Range("C3").Value = 19
$ws.Range("N3").Value = 23

$ws.Range("C4").Value = 28
$ws.Range("N4").Value = 72

$ws.Range("C5").Value = 11235.07
$ws.Range("N5").Value = 25567.58

$ws.Range("C6").Value = 15.04
$ws.Range("N6").Value = 34.24

$ws.Range("C7").Value = 401.25
$ws.Range("N7").Value = 355.11

$ws.Range("C8").Value = 132.16
$ws.Range("N8").Value = 74.22
